$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header columns G (latitude) and H (longitude), matching header style of existing headers
$ws.Range("A1").Copy($ws.Range("G1"))
$ws.Range("A1").Copy($ws.Range("H1"))
$ws.Cells.Item(1,7).Value = "latitude"
$ws.Cells.Item(1,8).Value = "longitude"

# Clear cells that become empty in the target data
$ws.Cells.Item(10,3).Value = $null
$ws.Cells.Item(14,3).Value = $null
$ws.Cells.Item(20,4).Value = $null

# Replace data rows 2-26 with the new dataset (A:H)
# row 2
$ws.Cells.Item(2,1).Value = 'Restaurant Y Panaderia Latina'
$ws.Cells.Item(2,2).Value = '3221 W Davis St, Dallas, TX 75211, États-Unis'
$ws.Cells.Item(2,3).Value = 'latinarestaurante.com'
$ws.Cells.Item(2,4).Value = '+1 214-337-4470'
$ws.Cells.Item(2,5).Value = 356
$ws.Cells.Item(2,6).Value = 4
$ws.Cells.Item(2,7).Value = 32.7498452
$ws.Cells.Item(2,8).Value = -97.1703813
# row 3
$ws.Cells.Item(3,1).Value = 'Sabor Latino'
$ws.Cells.Item(3,2).Value = '5431 E Grand Ave, Dallas, TX 75223, États-Unis'
$ws.Cells.Item(3,3).Value = 'saborlatinodallas.com'
$ws.Cells.Item(3,4).Value = '+1 817-243-5500'
$ws.Cells.Item(3,5).Value = 174
$ws.Cells.Item(3,6).Value = 4.3
$ws.Cells.Item(3,7).Value = 32.7498452
$ws.Cells.Item(3,8).Value = -97.1703813
# row 4
$ws.Cells.Item(4,1).Value = 'Miriam Cocina Latina'
$ws.Cells.Item(4,2).Value = '2015 Woodall Rodgers Fwy, Dallas, TX 75201, États-Unis'
$ws.Cells.Item(4,3).Value = 'miriamcocina.com'
$ws.Cells.Item(4,4).Value = '+1 214-855-5275'
$ws.Cells.Item(4,5).Value = 910
$ws.Cells.Item(4,6).Value = 4.6
$ws.Cells.Item(4,7).Value = 32.7498452
$ws.Cells.Item(4,8).Value = -97.1703813
# row 5
$ws.Cells.Item(5,1).Value = 'Zaguan Latin Café & Bakery'
$ws.Cells.Item(5,2).Value = 'Zaguan Latin Cafe & Bakery, 2604 Oak Lawn Ave, Dallas, TX 75219, États-Unis'
$ws.Cells.Item(5,3).Value = 'zaguan.com'
$ws.Cells.Item(5,4).Value = '+1 214-219-8393'
$ws.Cells.Item(5,5).Value = 1
$ws.Cells.Item(5,6).Value = 4.1
$ws.Cells.Item(5,7).Value = 32.7498452
$ws.Cells.Item(5,8).Value = -97.1703813
# row 6
$ws.Cells.Item(6,1).Value = 'Gloria''s Latin Cuisine'
$ws.Cells.Item(6,2).Value = '3715 Greenville Ave, Dallas, TX 75206, États-Unis'
$ws.Cells.Item(6,3).Value = 'gloriascuisine.com'
$ws.Cells.Item(6,4).Value = '+1 214-874-0088'
$ws.Cells.Item(6,5).Value = 1
$ws.Cells.Item(6,6).Value = 4.5
$ws.Cells.Item(6,7).Value = 32.7498452
$ws.Cells.Item(6,8).Value = -97.1703813
# row 7
$ws.Cells.Item(7,1).Value = 'Gloria''s Latin Cuisine'
$ws.Cells.Item(7,2).Value = '600 N Bishop Ave, Dallas, TX 75208, États-Unis'
$ws.Cells.Item(7,3).Value = 'gloriascuisine.com'
$ws.Cells.Item(7,4).Value = '+1 214-948-3672'
$ws.Cells.Item(7,5).Value = 2
$ws.Cells.Item(7,6).Value = 4.4
$ws.Cells.Item(7,7).Value = 32.7498452
$ws.Cells.Item(7,8).Value = -97.1703813
# row 8
$ws.Cells.Item(8,1).Value = 'Tienda & Restaurant Latina'
$ws.Cells.Item(8,2).Value = '4950 W Illinois Ave #115, Dallas, TX 75211, États-Unis'
$ws.Cells.Item(8,3).Value = 'latinarestaurante.com'
$ws.Cells.Item(8,4).Value = '+1 214-330-0011'
$ws.Cells.Item(8,5).Value = 272
$ws.Cells.Item(8,6).Value = 3.8
$ws.Cells.Item(8,7).Value = 32.7498452
$ws.Cells.Item(8,8).Value = -97.1703813
# row 9
$ws.Cells.Item(9,1).Value = 'Wild Salsa'
$ws.Cells.Item(9,2).Value = '1800 Main St, Dallas, TX 75201, États-Unis'
$ws.Cells.Item(9,3).Value = 'wildsalsarestaurant.com'
$ws.Cells.Item(9,4).Value = '+1 214-741-9453'
$ws.Cells.Item(9,5).Value = 1
$ws.Cells.Item(9,6).Value = 4.4
$ws.Cells.Item(9,7).Value = 32.7498452
$ws.Cells.Item(9,8).Value = -97.1703813
# row 10
$ws.Cells.Item(10,1).Value = 'Restaurant Latino'
$ws.Cells.Item(10,2).Value = '14500 S Josey Ln #101, Dallas, TX 75234, États-Unis'
$ws.Cells.Item(10,3).Value = $null
$ws.Cells.Item(10,4).Value = '+1 972-620-4954'
$ws.Cells.Item(10,5).Value = 172
$ws.Cells.Item(10,6).Value = 3.8
$ws.Cells.Item(10,7).Value = 32.9443412
$ws.Cells.Item(10,8).Value = -97.1783298
# row 11
$ws.Cells.Item(11,1).Value = 'Cafe Salsera'
$ws.Cells.Item(11,2).Value = '2610 Elm St, Dallas, TX 75226, États-Unis'
$ws.Cells.Item(11,3).Value = 'cafesalsera.com'
$ws.Cells.Item(11,4).Value = '+1 469-518-1500'
$ws.Cells.Item(11,5).Value = 682
$ws.Cells.Item(11,6).Value = 4
$ws.Cells.Item(11,7).Value = 32.9443412
$ws.Cells.Item(11,8).Value = -97.1783298
# row 12
$ws.Cells.Item(12,1).Value = 'El Bolero Cocina Mexicana'
$ws.Cells.Item(12,2).Value = '1201 Oak Lawn Ave #160, Dallas, TX 75207, États-Unis'
$ws.Cells.Item(12,3).Value = 'elboleromexican.com'
$ws.Cells.Item(12,4).Value = '+1 214-741-1986'
$ws.Cells.Item(12,5).Value = 1
$ws.Cells.Item(12,6).Value = 4.1
$ws.Cells.Item(12,7).Value = 32.9443412
$ws.Cells.Item(12,8).Value = -97.1783298
# row 13
$ws.Cells.Item(13,1).Value = 'Salsa Limón'
$ws.Cells.Item(13,2).Value = '411 N Akard St, Dallas, TX 75201, États-Unis'
$ws.Cells.Item(13,3).Value = 'salsalimon.com'
$ws.Cells.Item(13,4).Value = '+1 972-803-4388'
$ws.Cells.Item(13,5).Value = 1
$ws.Cells.Item(13,6).Value = 4.5
$ws.Cells.Item(13,7).Value = 32.9443412
$ws.Cells.Item(13,8).Value = -97.1783298
# row 14
$ws.Cells.Item(14,1).Value = 'Havana Cafe-Mojito Bar'
$ws.Cells.Item(14,2).Value = '1152 N Buckner Blvd #126, Dallas, TX 75218, États-Unis'
$ws.Cells.Item(14,3).Value = $null
$ws.Cells.Item(14,4).Value = '+1 214-660-9611'
$ws.Cells.Item(14,5).Value = 1
$ws.Cells.Item(14,6).Value = 4.4
$ws.Cells.Item(14,7).Value = 32.8342759
$ws.Cells.Item(14,8).Value = -96.9905284
# row 15
$ws.Cells.Item(15,1).Value = 'Cafe Madrid'
$ws.Cells.Item(15,2).Value = '4501 Travis St, Dallas, TX 75205, États-Unis'
$ws.Cells.Item(15,3).Value = 'cafemadrid-dallas.com'
$ws.Cells.Item(15,4).Value = '+1 214-528-1731'
$ws.Cells.Item(15,5).Value = 927
$ws.Cells.Item(15,6).Value = 4.4
$ws.Cells.Item(15,7).Value = 32.8342759
$ws.Cells.Item(15,8).Value = -96.9905284
# row 16
$ws.Cells.Item(16,1).Value = 'Gloria''s Latin Cuisine'
$ws.Cells.Item(16,2).Value = '3223 Lemmon Ave, Dallas, TX 75204, États-Unis'
$ws.Cells.Item(16,3).Value = 'gloriascuisine.com'
$ws.Cells.Item(16,4).Value = '+1 214-303-1166'
$ws.Cells.Item(16,5).Value = 1
$ws.Cells.Item(16,6).Value = 4.4
$ws.Cells.Item(16,7).Value = 32.8342759
$ws.Cells.Item(16,8).Value = -96.9905284
# row 17
$ws.Cells.Item(17,1).Value = 'Gloria''s Latin Cuisine'
$ws.Cells.Item(17,2).Value = '4140 Lemmon Ave, Dallas, TX 75219, États-Unis'
$ws.Cells.Item(17,3).Value = 'gloriascuisine.com'
$ws.Cells.Item(17,4).Value = '+1 214-521-7576'
$ws.Cells.Item(17,5).Value = 849
$ws.Cells.Item(17,6).Value = 4.4
$ws.Cells.Item(17,7).Value = 32.8342759
$ws.Cells.Item(17,8).Value = -96.9905284
# row 18
$ws.Cells.Item(18,1).Value = 'Las Palmas Restaurant'
$ws.Cells.Item(18,2).Value = '3957 Belt Line Rd, Addison, TX 75001, États-Unis'
$ws.Cells.Item(18,3).Value = 'laspalmasrestaurante.com'
$ws.Cells.Item(18,4).Value = '+1 469-665-8958'
$ws.Cells.Item(18,5).Value = 1
$ws.Cells.Item(18,6).Value = 4.1
$ws.Cells.Item(18,7).Value = 32.9539978
$ws.Cells.Item(18,8).Value = -97.1360389
# row 19
$ws.Cells.Item(19,1).Value = 'Beto & Son'
$ws.Cells.Item(19,2).Value = '3011 Gulden Ln #108, Dallas, TX 75212, États-Unis'
$ws.Cells.Item(19,3).Value = 'betoandsondallas.com'
$ws.Cells.Item(19,4).Value = '+1 469-249-8590'
$ws.Cells.Item(19,5).Value = 3
$ws.Cells.Item(19,6).Value = 4.3
$ws.Cells.Item(19,7).Value = 32.9539978
$ws.Cells.Item(19,8).Value = -97.1360389
# row 20
$ws.Cells.Item(20,1).Value = 'Ocho Latin Cuisine Events'
$ws.Cells.Item(20,2).Value = '369 Jefferson Blvd, Dallas, TX 75208, États-Unis'
$ws.Cells.Item(20,3).Value = 'ochooc.com'
$ws.Cells.Item(20,4).Value = $null
$ws.Cells.Item(20,5).Value = 134
$ws.Cells.Item(20,6).Value = 4.1
$ws.Cells.Item(20,7).Value = 32.7433967
$ws.Cells.Item(20,8).Value = -97.1164245
# row 21
$ws.Cells.Item(21,1).Value = 'Meso Maya Comida y Copas'
$ws.Cells.Item(21,2).Value = '1611 McKinney Ave, Dallas, TX 75202, États-Unis'
$ws.Cells.Item(21,3).Value = 'mesomaya.com'
$ws.Cells.Item(21,4).Value = '+1 214-484-6555'
$ws.Cells.Item(21,5).Value = 4
$ws.Cells.Item(21,6).Value = 4.5
$ws.Cells.Item(21,7).Value = 32.7433967
$ws.Cells.Item(21,8).Value = -97.1164245
# row 22
$ws.Cells.Item(22,1).Value = 'Campuzano Mexican Food'
$ws.Cells.Item(22,2).Value = '2618 Oak Lawn Ave, Dallas, TX 75219, États-Unis'
$ws.Cells.Item(22,3).Value = 'campuzanomexicanfood.com'
$ws.Cells.Item(22,4).Value = '+1 214-526-0100'
$ws.Cells.Item(22,5).Value = 2
$ws.Cells.Item(22,6).Value = 4.4
$ws.Cells.Item(22,7).Value = 32.7433967
$ws.Cells.Item(22,8).Value = -97.1164245
# row 23
$ws.Cells.Item(23,1).Value = 'La Duni Latin Cafe'
$ws.Cells.Item(23,2).Value = '2612 Irving Blvd, Dallas, TX 75207, États-Unis'
$ws.Cells.Item(23,3).Value = 'ladunihub.com'
$ws.Cells.Item(23,4).Value = '+1 214-520-7300'
$ws.Cells.Item(23,5).Value = 859
$ws.Cells.Item(23,6).Value = 4.3
$ws.Cells.Item(23,7).Value = 32.7433967
$ws.Cells.Item(23,8).Value = -97.1164245
# row 24
$ws.Cells.Item(24,1).Value = 'Mi Sazon Mexican Restaurant'
$ws.Cells.Item(24,2).Value = '3505 S Polk St, Dallas, TX 75224, États-Unis'
$ws.Cells.Item(24,3).Value = 'misazonrestaurant.com'
$ws.Cells.Item(24,4).Value = '+1 214-375-3333'
$ws.Cells.Item(24,5).Value = 745
$ws.Cells.Item(24,6).Value = 4.1
$ws.Cells.Item(24,7).Value = 32.7433967
$ws.Cells.Item(24,8).Value = -97.1164245
# row 25
$ws.Cells.Item(25,1).Value = 'Si Tapas'
$ws.Cells.Item(25,2).Value = '2207 Allen St, Dallas, TX 75204, États-Unis'
$ws.Cells.Item(25,3).Value = 'sitapasdallas.com'
$ws.Cells.Item(25,4).Value = '+1 214-720-0324'
$ws.Cells.Item(25,5).Value = 980
$ws.Cells.Item(25,6).Value = 4.3
$ws.Cells.Item(25,7).Value = 32.7433967
$ws.Cells.Item(25,8).Value = -97.1164245
# row 26
$ws.Cells.Item(26,1).Value = 'Te Deseo'
$ws.Cells.Item(26,2).Value = '2700 Olive St, Dallas, TX 75201, États-Unis'
$ws.Cells.Item(26,3).Value = 'tedeseo.com'
$ws.Cells.Item(26,4).Value = '+1 214-646-1314'
$ws.Cells.Item(26,5).Value = 1
$ws.Cells.Item(26,6).Value = 3.9
$ws.Cells.Item(26,7).Value = 32.7433967
$ws.Cells.Item(26,8).Value = -97.1164245

Write-Host "Applied google_maps_data update: added latitude/longitude columns and replaced sample rows."
